$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 545186.9399999999
$ws.Range("I9").Value = 669006.2
$ws.Range("J9").Value = 382.4
$ws.Range("K9").Value = 669006.2
$ws.Range("L9").Value = 382.4
$ws.Range("M9").Value = -668837.2
$ws.Range("N9").Value = -720.4

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2063.84
$ws.Range("J17").Value = 1785.0476
$ws.Range("L17").Value = 5355.142800000001
$ws.Range("N17").Value = -5691.142800000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 225.83333
$ws.Range("I31").Value = 225.83333
$ws.Range("K31").Value = 677.49999
$ws.Range("M31").Value = -447.49999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 26162.5
$ws.Range("I38").Value = 26162.5
$ws.Range("K38").Value = 78487.5
$ws.Range("M38").Value = -78115.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H42").Value = 334.5
$ws.Range("I42").Value = 169
$ws.Range("J42").Value = 500
$ws.Range("K42").Value = 507
$ws.Range("L42").Value = 1500
$ws.Range("M42").Value = -277
$ws.Range("N42").Value = -1960

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 2679.375
$ws.Range("I62").Value = 2355.8333
$ws.Range("J62").Value = 3650
$ws.Range("K62").Value = 2355.8333
$ws.Range("L62").Value = 3650
$ws.Range("M62").Value = -1731.8333
$ws.Range("N62").Value = -4898

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 90914740
$ws.Range("I64").Value = 5988.8887
$ws.Range("K64").Value = 5988.8887
$ws.Range("M64").Value = -5740.8887

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 2679.375
$ws.Range("I65").Value = 2355.8333
$ws.Range("J65").Value = 3650
$ws.Range("K65").Value = 11779.1665
$ws.Range("L65").Value = 18250
$ws.Range("M65").Value = -8659.166499999999
$ws.Range("N65").Value = -24490

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 90914740
$ws.Range("I67").Value = 5988.8887
$ws.Range("K67").Value = 5988.8887
$ws.Range("M67").Value = -5130.8887

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 5120.6
$ws.Range("I76").Value = 4901
$ws.Range("K76").Value = 4901
$ws.Range("M76").Value = -4586

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 5120.6
$ws.Range("I79").Value = 4901
$ws.Range("K79").Value = 4901
$ws.Range("M79").Value = -3809

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 5891519
$ws.Range("J88").Value = 12165.777
$ws.Range("L88").Value = 12165.777
$ws.Range("N88").Value = -12977.777

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H91").Value = 5891519
$ws.Range("J91").Value = 12165.777
$ws.Range("L91").Value = 12165.777
$ws.Range("N91").Value = -14973.777

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 4966.6665
$ws.Range("I98").Value = 4462.5
$ws.Range("K98").Value = 4462.5
$ws.Range("M98").Value = -2964.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 4966.6665
$ws.Range("I122").Value = 4462.5
$ws.Range("K122").Value = 13387.5
$ws.Range("M122").Value = -10937.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 3298452.2
$ws.Range("I132").Value = 6512.4614
$ws.Range("K132").Value = 19537.3842
$ws.Range("M132").Value = -17007.3842

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1788639.9
$ws.Range("I137").Value = 2383191.5
$ws.Range("J137").Value = 4985.4287
$ws.Range("K137").Value = 7149574.5
$ws.Range("L137").Value = 14956.2861
$ws.Range("M137").Value = -7147024.5
$ws.Range("N137").Value = -20056.2861

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 1500.1875
$ws.Range("I141").Value = 1419.5483
$ws.Range("K141").Value = 4258.644899999999
$ws.Range("M141").Value = 921.3551000000007

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2566.3635
$ws.Range("I32").Value = 2561.8147
$ws.Range("K32").Value = 2561.8147
$ws.Range("M32").Value = -2274.8147

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 11398516
$ws.Range("I132").Value = 2566021.5
$ws.Range("J132").Value = 55560990
$ws.Range("K132").Value = 7698064.5
$ws.Range("L132").Value = 166682970
$ws.Range("M132").Value = -7695534.5
$ws.Range("N132").Value = -166688030

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 162.47368
$ws.Range("I7").Value = 167.11111
$ws.Range("J7").Value = 79
$ws.Range("K7").Value = 167.11111
$ws.Range("L7").Value = 79
$ws.Range("M7").Value = -54.11111
$ws.Range("N7").Value = -305

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H33").Value = 749.5
$ws.Range("J33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("N33").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3739.6
$ws.Range("I58").Value = 3299.2856
$ws.Range("K58").Value = 3299.2856
$ws.Range("M58").Value = -3096.2856

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 1061.9524
$ws.Range("I94").Value = 628.875
$ws.Range("J94").Value = 1328.4615
$ws.Range("K94").Value = 628.875
$ws.Range("L94").Value = 1328.4615
$ws.Range("M94").Value = -177.875
$ws.Range("N94").Value = -2230.4615

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 2941961.8
$ws.Range("I107").Value = 5556098.5
$ws.Range("J107").Value = 1058.125
$ws.Range("K107").Value = 5556098.5
$ws.Range("L107").Value = 1058.125
$ws.Range("M107").Value = -5554178.5
$ws.Range("N107").Value = -4898.125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 3739.6
$ws.Range("I136").Value = 3299.2856
$ws.Range("K136").Value = 9897.856800000001
$ws.Range("M136").Value = -7347.856800000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H82").Value = 11500
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 11500
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 34500
$ws.Range("M82").ClearContents()
$ws.Range("N82").Value = -35312

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H85").Value = 11500
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 11500
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 34500
$ws.Range("M85").ClearContents()
$ws.Range("N85").Value = -37308

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 0
$ws.Range("I87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("M87").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H90").Value = 0
$ws.Range("I90").Value = 0
$ws.Range("K90").Value = 0
$ws.Range("M90").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H140").Value = 68792.914
$ws.Range("J140").Value = 68792.914
$ws.Range("L140").Value = 68792.914
$ws.Range("N140").Value = -79152.914

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("N40").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 279.72726
$ws.Range("I55").Value = 286.2
$ws.Range("K55").Value = 286.2
$ws.Range("M55").Value = -113.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H135").Value = 105995.5
$ws.Range("J135").Value = 105995.5
$ws.Range("L135").Value = 105995.5
$ws.Range("N135").Value = -116135.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H140").Value = 161549.72
$ws.Range("J140").Value = 161549.72
$ws.Range("L140").Value = 161549.72
$ws.Range("N140").Value = -171909.72

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 6999.5
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 6999.5
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 6999.5
$ws.Range("M96").ClearContents()
$ws.Range("N96").Value = -9745.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 5205.323
$ws.Range("I136").Value = 5487.9434
$ws.Range("J136").Value = 3957.0833
$ws.Range("K136").Value = 16463.8302
$ws.Range("L136").Value = 11871.2499
$ws.Range("M136").Value = -13913.8302
$ws.Range("N136").Value = -16971.2499

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H141").Value = 69998.73
$ws.Range("J141").Value = 69998.73
$ws.Range("L141").Value = 69998.73
$ws.Range("N141").Value = -80358.73
